# ---------------------------------------------------------------------------
# Applies the "added correct way to perform experiments to spreadsheet"
# commit to experiments.xlsx.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-FormatOnly {
    param([string]$From, [string]$To)
    $ws.Range($From).Copy() | Out-Null
    $ws.Range($To).PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Row 19 - add a "notes" comment in column K
# ---------------------------------------------------------------------------
$ws.Range("K19").Value = "Here there are only 2 loops. Strange because I did not perform the full unroll of the innermost loop."
Copy-FormatOnly "K18" "K19"

# ---------------------------------------------------------------------------
# Row 20 - only the row height changes (no content change)
# ---------------------------------------------------------------------------
$ws.Rows.Item(20).RowHeight = 133

# ---------------------------------------------------------------------------
# Row 21 - affine-loop-tile experiment
# ---------------------------------------------------------------------------
$ws.Range("D21").Value = "docker run -u `$(id -u) -v `$(pwd):/working_dir --rm agostini01/soda \`
                                                soda-opt \`
                                                -soda-outline-bambu-code \`
                                                -soda-extract-arguments-to-xml=using-bare-ptr \`
                                                -soda-generate-bambu-accelcode=no-aa \`
                                                -convert-linalg-to-affine-loops \`
                                                -affine-loop-tile \`
                                                -lower-all-to-llvm=use-bare-ptr-memref-call-conv \`
                                                -mlir-print-ir-after-all \`
                                                output/01searched-edited.mlir \`
                                                -o output/04optimized.mlir \`
                                                2>&1 | cat > output/05intermediate-optimized.mlir"
$ws.Range("C21").Value = "test to see if affine_loop tile has some impact on performances."
$ws.Range("E21").Value = "success"
Copy-FormatOnly "B21" "C21:E21"

$ws.Range("I21").Value = 33392
$ws.Range("J21").Formula = "=1-I21/H21"
Copy-FormatOnly "J19" "J21"
$ws.Rows.Item(21).RowHeight = 123

# ---------------------------------------------------------------------------
# Row 22 - affine-data-copy-generate experiment (failed)
# ---------------------------------------------------------------------------
$ws.Range("D22").Value = " docker run -u `$(id -u) -v `$(pwd):/working_dir --rm agostini01/soda \`
                                                  soda-opt \`
                                                  -soda-outline-bambu-code \`
                                                  -soda-extract-arguments-to-xml=using-bare-ptr \`
                                                  -soda-generate-bambu-accelcode=no-aa \`
                                                  -convert-linalg-to-affine-loops \`
                                                  -affine-data-copy-generate \`
                                                  -lower-all-to-llvm=use-bare-ptr-memref-call-conv \`
                                                  -mlir-print-ir-after-all \`
                                                  output/01searched-edited.mlir \`
                                                  -o output/04optimized.mlir \`
                                                  2>&1 | cat > output/05intermediate-optimized.mlir"
$ws.Range("C22").Value = "test to see if affine_data_copy_generate has some impact on performances."
$ws.Range("E22").Value = "error"
Copy-FormatOnly "B22" "C22:E22"

$ws.Range("I22").Value = "-"
$ws.Range("J22").Value = "-"
$ws.Rows.Item(22).RowHeight = 124

# ---------------------------------------------------------------------------
# Row 23 - affine-loop-unroll factor=1 experiment
# ---------------------------------------------------------------------------
$ws.Range("D23").Value = "docker run -u `$(id -u) -v `$(pwd):/working_dir --rm agostini01/soda \`
                                                soda-opt \`
                                                -soda-outline-bambu-code \`
                                                -soda-extract-arguments-to-xml=using-bare-ptr \`
                                                -soda-generate-bambu-accelcode \`
                                                -convert-linalg-to-affine-loops \`
                                                -affine-loop-unroll=`"unroll-factor=1`" \`
                                                -lower-all-to-llvm=use-bare-ptr-memref-call-conv \`
                                                -mlir-print-ir-after-all \`
                                                output/01searched-edited.mlir \`
                                                -o output/04optimized.mlir \`
                                                2>&1 | cat > output/05intermediate-optimized.mlir"
$ws.Range("C23").Value = "This is the right command to unroll the loop as desired (in combination with changes made to Bambu command). I am going to unroll more and more and analyze memory and computation bottlenecks. The unroll 1 has no effects on the loop, so it can be considered as a sort of baseline. Loop are slightly slower than baseline probably due to other some small optimizations in lower levels."
$ws.Range("E23").Value = "success"
Copy-FormatOnly "B23" "C23:E23"

$ws.Range("I23").Value = 29792
$ws.Range("J23").Formula = "=1-I23/H23"
Copy-FormatOnly "J19" "J23"
$ws.Rows.Item(23).RowHeight = 124

# ---------------------------------------------------------------------------
# Row 24 - affine-loop-unroll factor=2 experiment
# ---------------------------------------------------------------------------
$ws.Range("D24").Value = "docker run -u `$(id -u) -v `$(pwd):/working_dir --rm agostini01/soda \`
                                                soda-opt \`
                                                -soda-outline-bambu-code \`
                                                -soda-extract-arguments-to-xml=using-bare-ptr \`
                                                -soda-generate-bambu-accelcode \`
                                                -convert-linalg-to-affine-loops \`
                                                -affine-loop-unroll=`"unroll-factor=2`" \`
                                                -lower-all-to-llvm=use-bare-ptr-memref-call-conv \`
                                                -mlir-print-ir-after-all \`
                                                output/01searched-edited.mlir \`
                                                -o output/04optimized.mlir \`
                                                2>&1 | cat > output/05intermediate-optimized.mlir"
$ws.Range("C24").Value = "This is the right command to unroll the loop as desired (in combination with changes made to Bambu command). I am going to unroll more and more and analyze memory and computation bottlenecks."
$ws.Range("E24").Value = "success"
Copy-FormatOnly "B24" "C24:E24"

$ws.Range("I24").Value = 23117
$ws.Range("J24").Formula = "=1-I24/H24"
Copy-FormatOnly "J19" "J24"
$ws.Rows.Item(24).RowHeight = 120

# ---------------------------------------------------------------------------
# Row 25 - affine-loop-unroll factor=3 experiment
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = "docker run -u `$(id -u) -v `$(pwd):/working_dir --rm agostini01/soda \`
                                                soda-opt \`
                                                -soda-outline-bambu-code \`
                                                -soda-extract-arguments-to-xml=using-bare-ptr \`
                                                -soda-generate-bambu-accelcode \`
                                                -convert-linalg-to-affine-loops \`
                                                -affine-loop-unroll=`"unroll-factor=3`" \`
                                                -lower-all-to-llvm=use-bare-ptr-memref-call-conv \`
                                                -mlir-print-ir-after-all \`
                                                output/01searched-edited.mlir \`
                                                -o output/04optimized.mlir \`
                                                2>&1 | cat > output/05intermediate-optimized.mlir"
$ws.Range("C25").Value = "In this case the unrolling behaves as expected, making three unrolls of the loop. In this way the innermost loop has no more 15 iterations but 5 (15 step 3)."
Copy-FormatOnly "B25" "C25:D25"

$ws.Range("I25").Value = 20192
$ws.Range("J25").Formula = "=1-I25/H25"
Copy-FormatOnly "J19" "J25"
$ws.Rows.Item(25).RowHeight = 112

# ---------------------------------------------------------------------------
# Row 29 - new "matmul 15" baseline row (was a blank filler row)
# ---------------------------------------------------------------------------
$ws.Range("D29").ClearContents() | Out-Null
$ws.Range("D29").ClearFormats() | Out-Null

$ws.Range("B29").Value = "matmul 15"
Copy-FormatOnly "B21" "B29"

$ws.Range("F29").Value = 24000
$ws.Range("H29").Value = 33392
Copy-FormatOnly "F21:I21" "F29:I29"

Copy-FormatOnly "J26" "J29"

$ws.Rows.Item(29).RowHeight = 17

# ---------------------------------------------------------------------------
# Rows 39, 41, 42 - extend the trailing blank-formatted filler rows by one
# ---------------------------------------------------------------------------
Copy-FormatOnly "D30" "D39"
Copy-FormatOnly "F40:I40" "F41:I41"
Copy-FormatOnly "J41" "J42"

# ---------------------------------------------------------------------------
# Sheet view - scroll position + active selection
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("J25").Select() | Out-Null

Write-Output "edit complete"
